$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the topic text for row 10 (Deployment + Runtime -> Deployment + Runtime + Design Decisions)
$ws.Range("B10").Value = "Deployment + Runtime + Design Decisions"

# Update effort value for row 10 (3 -> 4 hrs), which also recalculates the Total effort formula in C11
$ws.Range("C10").Value = 4

# Row 10 height changes to 30
$ws.Rows.Item(10).RowHeight = 30

# Update the sheet view: move selection to H25
$ws.Range("H25").Select()
